{"js": "// Replace the date header and each two-digit multiplication problem's\n// text with its updated value. Every \"old\" value below is unique within\n// the document, so a body-wide search for each one is unambiguous.\nconst replacements = [\n  [\"2023-03-08 Wednesday\", \"2023-03-09 Thursday\"],\n  [\"45\u00d795=\", \"50\u00d717=\"],\n  [\"58\u00d750=\", \"59\u00d777=\"],\n  [\"65\u00d731=\", \"97\u00d729=\"],\n  [\"24\u00d777=\", \"16\u00d793=\"],\n  [\"98\u00d752=\", \"50\u00d765=\"],\n  [\"37\u00d733=\", \"80\u00d768=\"],\n  [\"67\u00d731=\", \"98\u00d723=\"],\n  [\"11\u00d747=\", \"98\u00d771=\"],\n  [\"75\u00d735=\", \"54\u00d797=\"],\n  [\"41\u00d789=\", \"80\u00d718=\"],\n  [\"23\u00d711=\", \"58\u00d724=\"],\n  [\"63\u00d764=\", \"56\u00d715=\"],\n  [\"95\u00d732=\", \"23\u00d7100=\"],\n  [\"97\u00d714=\", \"70\u00d712=\"],\n  [\"40\u00d788=\", \"54\u00d758=\"],\n  [\"72\u00d711=\", \"87\u00d761=\"],\n  [\"45\u00d733=\", \"55\u00d783=\"],\n  [\"29\u00d716=\", \"99\u00d712=\"],\n  [\"50\u00d723=\", \"13\u00d757=\"],\n  [\"60\u00d779=\", \"82\u00d752=\"],\n  [\"43\u00d764=\", \"51\u00d778=\"],\n  [\"18\u00d722=\", \"66\u00d786=\"],\n  [\"51\u00d759=\", \"21\u00d759=\"],\n  [\"93\u00d769=\", \"92\u00d731=\"],\n  [\"100\u00d770=\", \"89\u00d774=\"],\n  [\"33\u00d795=\", \"26\u00d756=\"],\n  [\"11\u00d721=\", \"50\u00d728=\"],\n  [\"90\u00d799=\", \"27\u00d711=\"],\n  [\"69\u00d754=\", \"21\u00d731=\"],\n  [\"38\u00d760=\", \"88\u00d786=\"],\n  [\"64\u00d727=\", \"30\u00d719=\"],\n  [\"75\u00d723=\", \"89\u00d727=\"],\n  [\"54\u00d756=\", \"98\u00d711=\"],\n  [\"85\u00d792=\", \"32\u00d723=\"],\n  [\"83\u00d755=\", \"14\u00d790=\"],\n  [\"63\u00d784=\", \"51\u00d741=\"],\n  [\"81\u00d733=\", \"57\u00d744=\"],\n  [\"78\u00d785=\", \"34\u00d768=\"],\n  [\"43\u00d796=\", \"34\u00d726=\"],\n  [\"86\u00d756=\", \"53\u00d786=\"],\n  [\"25\u00d762=\", \"84\u00d733=\"],\n  [\"58\u00d782=\", \"64\u00d745=\"],\n  [\"76\u00d717=\", \"88\u00d792=\"],\n  [\"81\u00d772=\", \"56\u00d714=\"],\n  [\"22\u00d754=\", \"79\u00d718=\"],\n  [\"81\u00d715=\", \"18\u00d798=\"],\n  [\"41\u00d746=\", \"72\u00d778=\"],\n  [\"35\u00d740=\", \"39\u00d774=\"],\n  [\"62\u00d714=\", \"43\u00d759=\"],\n  [\"97\u00d780=\", \"48\u00d741=\"],\n  [\"96\u00d718=\", \"82\u00d722=\"],\n  [\"27\u00d730=\", \"35\u00d720=\"],\n  [\"57\u00d771=\", \"36\u00d732=\"],\n  [\"43\u00d717=\", \"11\u00d743=\"],\n  [\"70\u00d782=\", \"29\u00d781=\"],\n  [\"94\u00d798=\", \"40\u00d741=\"],\n  [\"22\u00d719=\", \"69\u00d790=\"],\n  [\"84\u00d768=\", \"88\u00d758=\"],\n  [\"87\u00d758=\", \"39\u00d745=\"],\n  [\"19\u00d797=\", \"15\u00d729=\"],\n  [\"19\u00d788=\", \"92\u00d797=\"],\n  [\"16\u00d779=\", \"35\u00d719=\"],\n  [\"100\u00d767=\", \"85\u00d757=\"],\n  [\"61\u00d730=\", \"94\u00d763=\"],\n  [\"100\u00d764=\", \"75\u00d761=\"],\n  [\"22\u00d792=\", \"82\u00d789=\"],\n  [\"14\u00d730=\", \"51\u00d752=\"],\n  [\"15\u00d787=\", \"98\u00d784=\"],\n  [\"90\u00d747=\", \"42\u00d715=\"],\n  [\"51\u00d733=\", \"42\u00d724=\"],\n  [\"15\u00d756=\", \"40\u00d714=\"],\n  [\"75\u00d725=\", \"28\u00d751=\"],\n  [\"12\u00d766=\", \"89\u00d777=\"],\n  [\"98\u00d767=\", \"40\u00d783=\"],\n  [\"28\u00d796=\", \"46\u00d762=\"],\n  [\"81\u00d739=\", \"90\u00d722=\"],\n  [\"51\u00d7100=\", \"30\u00d762=\"],\n  [\"46\u00d766=\", \"98\u00d754=\"],\n  [\"67\u00d785=\", \"22\u00d799=\"],\n  [\"52\u00d787=\", \"58\u00d770=\"],\n  [\"47\u00d783=\", \"81\u00d737=\"],\n  [\"86\u00d788=\", \"86\u00d746=\"],\n  [\"40\u00d735=\", \"45\u00d771=\"],\n  [\"65\u00d745=\", \"31\u00d765=\"],\n  [\"28\u00d799=\", \"64\u00d741=\"],\n  [\"25\u00d744=\", \"14\u00d715=\"],\n  [\"86\u00d781=\", \"16\u00d788=\"],\n  [\"71\u00d755=\", \"43\u00d765=\"],\n  [\"24\u00d760=\", \"83\u00d746=\"],\n  [\"33\u00d760=\", \"78\u00d789=\"],\n  [\"13\u00d718=\", \"95\u00d764=\"],\n  [\"16\u00d754=\", \"29\u00d739=\"],\n  [\"87\u00d727=\", \"35\u00d749=\"],\n  [\"91\u00d726=\", \"16\u00d756=\"],\n  [\"10\u00d731=\", \"99\u00d768=\"],\n  [\"32\u00d772=\", \"80\u00d761=\"],\n  [\"49\u00d736=\", \"84\u00d741=\"],\n  [\"12\u00d794=\", \"71\u00d781=\"],\n  [\"28\u00d745=\", \"29\u00d762=\"],\n  [\"60\u00d720=\", \"97\u00d781=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date header and each two-digit multiplication problem's\n# text with its updated value. Every \"old\" value below is unique within\n# the document, so Find/Replace for each one is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-03-08 Wednesday\", \"2023-03-09 Thursday\"),\n    @(\"45\u00d795=\", \"50\u00d717=\"),\n    @(\"58\u00d750=\", \"59\u00d777=\"),\n    @(\"65\u00d731=\", \"97\u00d729=\"),\n    @(\"24\u00d777=\", \"16\u00d793=\"),\n    @(\"98\u00d752=\", \"50\u00d765=\"),\n    @(\"37\u00d733=\", \"80\u00d768=\"),\n    @(\"67\u00d731=\", \"98\u00d723=\"),\n    @(\"11\u00d747=\", \"98\u00d771=\"),\n    @(\"75\u00d735=\", \"54\u00d797=\"),\n    @(\"41\u00d789=\", \"80\u00d718=\"),\n    @(\"23\u00d711=\", \"58\u00d724=\"),\n    @(\"63\u00d764=\", \"56\u00d715=\"),\n    @(\"95\u00d732=\", \"23\u00d7100=\"),\n    @(\"97\u00d714=\", \"70\u00d712=\"),\n    @(\"40\u00d788=\", \"54\u00d758=\"),\n    @(\"72\u00d711=\", \"87\u00d761=\"),\n    @(\"45\u00d733=\", \"55\u00d783=\"),\n    @(\"29\u00d716=\", \"99\u00d712=\"),\n    @(\"50\u00d723=\", \"13\u00d757=\"),\n    @(\"60\u00d779=\", \"82\u00d752=\"),\n    @(\"43\u00d764=\", \"51\u00d778=\"),\n    @(\"18\u00d722=\", \"66\u00d786=\"),\n    @(\"51\u00d759=\", \"21\u00d759=\"),\n    @(\"93\u00d769=\", \"92\u00d731=\"),\n    @(\"100\u00d770=\", \"89\u00d774=\"),\n    @(\"33\u00d795=\", \"26\u00d756=\"),\n    @(\"11\u00d721=\", \"50\u00d728=\"),\n    @(\"90\u00d799=\", \"27\u00d711=\"),\n    @(\"69\u00d754=\", \"21\u00d731=\"),\n    @(\"38\u00d760=\", \"88\u00d786=\"),\n    @(\"64\u00d727=\", \"30\u00d719=\"),\n    @(\"75\u00d723=\", \"89\u00d727=\"),\n    @(\"54\u00d756=\", \"98\u00d711=\"),\n    @(\"85\u00d792=\", \"32\u00d723=\"),\n    @(\"83\u00d755=\", \"14\u00d790=\"),\n    @(\"63\u00d784=\", \"51\u00d741=\"),\n    @(\"81\u00d733=\", \"57\u00d744=\"),\n    @(\"78\u00d785=\", \"34\u00d768=\"),\n    @(\"43\u00d796=\", \"34\u00d726=\"),\n    @(\"86\u00d756=\", \"53\u00d786=\"),\n    @(\"25\u00d762=\", \"84\u00d733=\"),\n    @(\"58\u00d782=\", \"64\u00d745=\"),\n    @(\"76\u00d717=\", \"88\u00d792=\"),\n    @(\"81\u00d772=\", \"56\u00d714=\"),\n    @(\"22\u00d754=\", \"79\u00d718=\"),\n    @(\"81\u00d715=\", \"18\u00d798=\"),\n    @(\"41\u00d746=\", \"72\u00d778=\"),\n    @(\"35\u00d740=\", \"39\u00d774=\"),\n    @(\"62\u00d714=\", \"43\u00d759=\"),\n    @(\"97\u00d780=\", \"48\u00d741=\"),\n    @(\"96\u00d718=\", \"82\u00d722=\"),\n    @(\"27\u00d730=\", \"35\u00d720=\"),\n    @(\"57\u00d771=\", \"36\u00d732=\"),\n    @(\"43\u00d717=\", \"11\u00d743=\"),\n    @(\"70\u00d782=\", \"29\u00d781=\"),\n    @(\"94\u00d798=\", \"40\u00d741=\"),\n    @(\"22\u00d719=\", \"69\u00d790=\"),\n    @(\"84\u00d768=\", \"88\u00d758=\"),\n    @(\"87\u00d758=\", \"39\u00d745=\"),\n    @(\"19\u00d797=\", \"15\u00d729=\"),\n    @(\"19\u00d788=\", \"92\u00d797=\"),\n    @(\"16\u00d779=\", \"35\u00d719=\"),\n    @(\"100\u00d767=\", \"85\u00d757=\"),\n    @(\"61\u00d730=\", \"94\u00d763=\"),\n    @(\"100\u00d764=\", \"75\u00d761=\"),\n    @(\"22\u00d792=\", \"82\u00d789=\"),\n    @(\"14\u00d730=\", \"51\u00d752=\"),\n    @(\"15\u00d787=\", \"98\u00d784=\"),\n    @(\"90\u00d747=\", \"42\u00d715=\"),\n    @(\"51\u00d733=\", \"42\u00d724=\"),\n    @(\"15\u00d756=\", \"40\u00d714=\"),\n    @(\"75\u00d725=\", \"28\u00d751=\"),\n    @(\"12\u00d766=\", \"89\u00d777=\"),\n    @(\"98\u00d767=\", \"40\u00d783=\"),\n    @(\"28\u00d796=\", \"46\u00d762=\"),\n    @(\"81\u00d739=\", \"90\u00d722=\"),\n    @(\"51\u00d7100=\", \"30\u00d762=\"),\n    @(\"46\u00d766=\", \"98\u00d754=\"),\n    @(\"67\u00d785=\", \"22\u00d799=\"),\n    @(\"52\u00d787=\", \"58\u00d770=\"),\n    @(\"47\u00d783=\", \"81\u00d737=\"),\n    @(\"86\u00d788=\", \"86\u00d746=\"),\n    @(\"40\u00d735=\", \"45\u00d771=\"),\n    @(\"65\u00d745=\", \"31\u00d765=\"),\n    @(\"28\u00d799=\", \"64\u00d741=\"),\n    @(\"25\u00d744=\", \"14\u00d715=\"),\n    @(\"86\u00d781=\", \"16\u00d788=\"),\n    @(\"71\u00d755=\", \"43\u00d765=\"),\n    @(\"24\u00d760=\", \"83\u00d746=\"),\n    @(\"33\u00d760=\", \"78\u00d789=\"),\n    @(\"13\u00d718=\", \"95\u00d764=\"),\n    @(\"16\u00d754=\", \"29\u00d739=\"),\n    @(\"87\u00d727=\", \"35\u00d749=\"),\n    @(\"91\u00d726=\", \"16\u00d756=\"),\n    @(\"10\u00d731=\", \"99\u00d768=\"),\n    @(\"32\u00d772=\", \"80\u00d761=\"),\n    @(\"49\u00d736=\", \"84\u00d741=\"),\n    @(\"12\u00d794=\", \"71\u00d781=\"),\n    @(\"28\u00d745=\", \"29\u00d762=\"),\n    @(\"60\u00d720=\", \"97\u00d781=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdReplaceAll=2: replace every occurrence of this (unique) old value\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
